# Update the "Corr/total marks" figures in the concise marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: Right (correct) marks total 3 -> 5
$ws.Range("B11").Value = 5

# Total row: Right (correct) marks total 15 -> 25
$ws.Range("B12").Value = 25

# Total row: Max column "corr/total" text 9/84 -> 25/140
$ws.Range("E12").Value = "25/140"
